$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeit")

# Translate the "Arbeit" job/animal names (B2:B15) from English to Korean.
$ws.Range("B2").Value  = "여우"
$ws.Range("B3").Value  = "학"
$ws.Range("B4").Value  = "사자"
$ws.Range("B5").Value  = "두꺼비"
$ws.Range("B6").Value  = "소"
$ws.Range("B7").Value  = "늑대"
$ws.Range("B8").Value  = "곰"
$ws.Range("B9").Value  = "개"
$ws.Range("B10").Value = "거위"
$ws.Range("B11").Value = "토끼"
$ws.Range("B12").Value = "쥐"
$ws.Range("B13").Value = "돼지"
$ws.Range("B14").Value = "제비"
$ws.Range("B15").Value = "박쥐"

# Fix which tab/cell is active: "Arbeit" should be the selected sheet (not
# "ArbeitHappiness"), with B15 selected instead of the previous H1 / L12.
$ws.Activate()
$ws.Range("B15").Select()
